$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12
$ws.Range("F12").Value = 1305498328.7699957
$ws.Range("G12").Value = 12405247717

# Row 13
$ws.Range("F13").Value = 325268233.58999997
$ws.Range("G13").Value = 319819483.19999999

# Row 14
$ws.Range("G14").Value = 34063116.799999997

# Row 16
$ws.Range("F16").Value = -53616441.74000001

# Row 18 (F18 becomes a formula; G18 keeps its existing formula and recalculates)
$ws.Range("F18").Formula = "=SUM(F12:F17)"

# Row 19
$ws.Range("F19").Value = -412700000
$ws.Range("G19").Value = -379300000

# Row 21 (F21 becomes a formula; G21 keeps its existing formula and recalculates)
$ws.Range("F21").Formula = "=SUM(F18:F20)"

# Row 22
$ws.Range("F22").Value = -20015625

# Row 26
$ws.Range("F26").Value = 1026703455.3810816
$ws.Range("G26").Value = 1029174575

# Force full recalculation so dependent formulas (F18 dep chain, G18, F21, G21,
# F23, G23, F25, G25, F28, G28, C29) pick up the new values.
$excel.CalculateFullRebuild()
